$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value updates (row 3) ---
$ws.Range("S3").Value = 0.82
$ws.Range("T3").Value = 0.82
$ws.Range("U3").Value = 0.82

# --- Cell value updates (row 5) ---
$ws.Range("G5").Value = 0.654
$ws.Range("M5").Value = 0.66

# --- Cell value updates (row 8) ---
$ws.Range("Q8").Value = 0.819

# --- Cell value updates (row 9) ---
$ws.Range("S9").Value = 0.819
$ws.Range("T9").Value = 0.819
$ws.Range("U9").Value = 0.819
$ws.Range("V9").Value = 0.819

# --- Cell value updates (row 10) ---
$ws.Range("S10").Value = 0.82
$ws.Range("T10").Value = 0.819
$ws.Range("U10").Value = 0.819
$ws.Range("V10").Value = 0.82
$ws.Range("W10").Value = 0.82

# --- Cell value updates (row 11) ---
$ws.Range("T11").Value = 0.819

# --- Cell value updates (row 15) ---
$ws.Range("W15").Value = 0.86

# --- Cell value updates (row 20) ---
$ws.Range("T20").Value = 0.859

# --- Column widths for B:N (explicit default-like width, ~9.140625) ---
$ws.Range("B1:N1").EntireColumn.ColumnWidth = 8.4

# --- View / selection changes ---
$ws.Activate()
$ws.Range("H25").Select()

# --- Best-effort: restore the window size/position (no visual effect on data) ---
try {
    $excel.ActiveWindow.Left = 1230
    $excel.ActiveWindow.Top = 2985
    $excel.ActiveWindow.Width = 27570
    $excel.ActiveWindow.Height = 12555
} catch {
}
